# Add a "Turkey" test-data sheet, cloned from "Spain", for the Zettler
# printer template (per commit: "Added Test data for Turkey Template for
# Zettler").

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# Clone "Spain" (keeps styles/merges/column widths/etc. identical) and place
# the copy immediately after it.
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item("Spain (2)")
$turkey.Name = "Turkey"

# Turkey-specific content (new shared strings "Turkey Market" /
# "NGC-3191/T3300 ").
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3300 "

# Column D needs to be a bit wider on the new sheet so the longer "Expected
# value" text no longer wraps.
$turkey.Columns.Item(4).ColumnWidth = 24.83
$turkey.Rows.Item(5).AutoFit()

# Restore Spain's own selection (it's no longer the active/selected tab)
# and leave Turkey as the freshly active tab with its own selection.
$spain.Range("A1:D15").Select()
$turkey.Activate()
$turkey.Range("N3").Select()
